$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 2763.0557
$ws.Range("I15").Value = 2763.0557
$ws.Range("K15").Value = 8289.167099999999
$ws.Range("M15").Value = -8120.167099999999

$ws.Range("H39").Value = 212.73685
$ws.Range("I39").Value = 120.166664
$ws.Range("K39").Value = 360.499992
$ws.Range("M39").Value = -64.49999200000002

$ws.Range("H43").Value = 3500
$ws.Range("J43").Value = 3500
$ws.Range("L43").Value = 3500
$ws.Range("N43").Value = -3638

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()

$ws.Range("H62").Value = 1916.625
$ws.Range("J62").Value = 2250
$ws.Range("L62").Value = 2250
$ws.Range("N62").Value = -3498

$ws.Range("H65").Value = 1916.625
$ws.Range("J65").Value = 2250
$ws.Range("L65").Value = 11250
$ws.Range("N65").Value = -17490

$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

$ws.Range("H128").Value = 38571.43
$ws.Range("J128").Value = 38571.43
$ws.Range("L128").Value = 38571.43
$ws.Range("N128").Value = -48531.43

$ws.Range("H132").Value = 13710.596
$ws.Range("I132").Value = 2769.1562
$ws.Range("K132").Value = 8307.4686
$ws.Range("M132").Value = -5777.4686

$ws.Range("H137").Value = 2024.3334
$ws.Range("I137").Value = 1877.6
$ws.Range("K137").Value = 5632.799999999999
$ws.Range("M137").Value = -3082.799999999999

$ws.Range("H138").Value = 3501.6938
$ws.Range("I138").Value = 2606.5386
$ws.Range("J138").Value = 3824.9443
$ws.Range("K138").Value = 7819.6158
$ws.Range("L138").Value = 11474.8329
$ws.Range("M138").Value = -2679.6158
$ws.Range("N138").Value = -21754.8329

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3260.0364
$ws.Range("I32").Value = 3314.4807
$ws.Range("K32").Value = 3314.4807
$ws.Range("M32").Value = -3027.4807

$ws.Range("H74").Value = 1351.0209
$ws.Range("I74").Value = 1428.3823
$ws.Range("K74").Value = 1428.3823
$ws.Range("M74").Value = -554.3823

$ws.Range("H77").Value = 1351.0209
$ws.Range("I77").Value = 1428.3823
$ws.Range("K77").Value = 7141.9115
$ws.Range("M77").Value = -2773.9115

$ws.Range("H132").Value = 1925.7858
$ws.Range("I132").Value = 1766.0952
$ws.Range("K132").Value = 5298.2856
$ws.Range("M132").Value = -2768.2856

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 17242824
$ws.Range("I86").Value = 20834760
$ws.Range("J86").Value = 1531.6
$ws.Range("K86").Value = 20834760
$ws.Range("L86").Value = 1531.6
$ws.Range("M86").Value = -20833637
$ws.Range("N86").Value = -3777.6

$ws.Range("H87").Value = 50000
$ws.Range("J87").Value = 50000
$ws.Range("L87").Value = 50000
$ws.Range("N87").Value = -52496

$ws.Range("H89").Value = 17242824
$ws.Range("I89").Value = 20834760
$ws.Range("J89").Value = 1531.6
$ws.Range("K89").Value = 104173800
$ws.Range("L89").Value = 7658
$ws.Range("M89").Value = -104168184
$ws.Range("N89").Value = -18890

$ws.Range("H90").Value = 50000
$ws.Range("J90").Value = 50000
$ws.Range("L90").Value = 150000
$ws.Range("N90").Value = -162480

$ws.Range("H107").Value = 7098.625
$ws.Range("I107").Value = 6081.923
$ws.Range("K107").Value = 6081.923
$ws.Range("M107").Value = -4161.923

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 322.5
$ws.Range("I7").Value = 157.5
$ws.Range("J7").Value = 432.5
$ws.Range("K7").Value = 157.5
$ws.Range("L7").Value = 432.5
$ws.Range("M7").Value = -44.5
$ws.Range("N7").Value = -658.5

$ws.Range("H58").Value = 2647.8333
$ws.Range("I58").Value = 1666.6666
$ws.Range("J58").Value = 2974.889
$ws.Range("K58").Value = 1666.6666
$ws.Range("L58").Value = 2974.889
$ws.Range("M58").Value = -1463.6666
$ws.Range("N58").Value = -3380.889

$ws.Range("H99").Value = 4405
$ws.Range("I99").Value = 5296.75
$ws.Range("J99").Value = 3216
$ws.Range("K99").Value = 5296.75
$ws.Range("L99").Value = 3216
$ws.Range("M99").Value = -3798.75
$ws.Range("N99").Value = -6212

$ws.Range("H122").Value = 4763.5454
$ws.Range("I122").Value = 4178.222
$ws.Range("J122").Value = 5168.769
$ws.Range("K122").Value = 12534.666
$ws.Range("L122").Value = 15506.307
$ws.Range("M122").Value = -10084.666
$ws.Range("N122").Value = -20406.307

$ws.Range("H126").Value = 4405
$ws.Range("I126").Value = 5296.75
$ws.Range("J126").Value = 3216
$ws.Range("K126").Value = 15890.25
$ws.Range("L126").Value = 9648
$ws.Range("M126").Value = -13420.25
$ws.Range("N126").Value = -14588

$ws.Range("H132").Value = 4880
$ws.Range("I132").Value = 3700.5
$ws.Range("J132").Value = 5666.3335
$ws.Range("K132").Value = 11101.5
$ws.Range("L132").Value = 16999.0005
$ws.Range("M132").Value = -8571.5
$ws.Range("N132").Value = -22059.0005

$ws.Range("H134").Value = 5561.5
$ws.Range("I134").Value = 5748.8335
$ws.Range("J134").Value = 4999.5
$ws.Range("K134").Value = 17246.5005
$ws.Range("L134").Value = 14998.5
$ws.Range("M134").Value = -14711.5005
$ws.Range("N134").Value = -20068.5

$ws.Range("H136").Value = 2647.8333
$ws.Range("I136").Value = 1666.6666
$ws.Range("J136").Value = 2974.889
$ws.Range("K136").Value = 4999.9998
$ws.Range("L136").Value = 8924.667000000001
$ws.Range("M136").Value = -2449.9998
$ws.Range("N136").Value = -14024.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 148015.12
$ws.Range("I70").Value = 191686.83
$ws.Range("J70").Value = 17000
$ws.Range("K70").Value = 191686.83
$ws.Range("L70").Value = 17000
$ws.Range("M70").Value = -191416.83
$ws.Range("N70").Value = -17540

$ws.Range("H73").Value = 148015.12
$ws.Range("I73").Value = 191686.83
$ws.Range("J73").Value = 17000
$ws.Range("K73").Value = 191686.83
$ws.Range("L73").Value = 17000
$ws.Range("M73").Value = -190750.83
$ws.Range("N73").Value = -18872

$ws.Range("H107").Value = 490.6842
$ws.Range("I107").Value = 366.4375
$ws.Range("K107").Value = 366.4375
$ws.Range("M107").Value = 1553.5625

$ws.Range("H122").Value = 9236.75
$ws.Range("I122").Value = 950
$ws.Range("J122").Value = 11999
$ws.Range("K122").Value = 2850
$ws.Range("L122").Value = 35997
$ws.Range("M122").Value = -400
$ws.Range("N122").Value = -40897

$ws.Range("H132").Value = 6660.59
$ws.Range("I132").Value = 5787.1763
$ws.Range("K132").Value = 17361.5289
$ws.Range("M132").Value = -14831.5289

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 312.7143
$ws.Range("I32").Value = 356.5
$ws.Range("K32").Value = 356.5
$ws.Range("M32").Value = -39.5

$ws.Range("H40").Value = 2902.4
$ws.Range("I40").Value = 2651.75
$ws.Range("K40").Value = 2651.75
$ws.Range("M40").Value = -2515.75

$ws.Range("H46").Value = 2615.12
$ws.Range("I46").Value = 1769.1
$ws.Range("J46").Value = 3179.1333
$ws.Range("K46").Value = 1769.1
$ws.Range("L46").Value = 3179.1333
$ws.Range("M46").Value = -1581.1
$ws.Range("N46").Value = -3555.1333

$ws.Range("H131").Value = 39750
$ws.Range("J131").Value = 39750
$ws.Range("L131").Value = 39750
$ws.Range("N131").Value = -49830

$ws.Range("H132").Value = 2755.9
$ws.Range("I132").Value = 1874.5333
$ws.Range("J132").Value = 5400
$ws.Range("K132").Value = 5623.5999
$ws.Range("L132").Value = 16200
$ws.Range("M132").Value = -3093.5999
$ws.Range("N132").Value = -21260

$ws.Range("H136").Value = 4383.8423
$ws.Range("I136").Value = 3949.3
$ws.Range("K136").Value = 11847.9
$ws.Range("M136").Value = -9297.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 570
$ws.Range("J113").Value = 708.5
$ws.Range("L113").Value = 2125.5
$ws.Range("N113").Value = -6465.5

$ws.Range("H122").Value = 1949.75
$ws.Range("I122").Value = 1949.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5849.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -3399.25
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 4733.1665
$ws.Range("I132").Value = 4099.75
$ws.Range("K132").Value = 12299.25
$ws.Range("M132").Value = -9769.25

$ws.Range("H136").Value = 1610.5555
$ws.Range("I136").Value = 1499.4286
$ws.Range("K136").Value = 4498.2858
$ws.Range("M136").Value = -1948.2858
